# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header formatting used by the other header cells (B1:G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring header cell (G1) onto the new
# header cell (H1) so it reuses the existing header style rather than
# creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for row 2 under the "Save" column.
$ws.Range("H2").Value = 0
